$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete the row containing "Essendon v Collingwood" (entire row, shifts rows below up)
$ws.Rows.Item(4).Delete()

# Update home_line (B) and away_line (C) values for the remaining matches
# Row 2: Brisbane Lions v St Kilda
$ws.Range("B2").Value = -29.25
$ws.Range("C2").Value = 29.25

# Row 3: Carlton v Greater Western Sydney
$ws.Range("B3").Value = -3.5
$ws.Range("C3").Value = 3.5

# Row 4 (was row 5): Geelong v Western Bulldogs
$ws.Range("B4").Value = 8.5
$ws.Range("C4").Value = -8.5

# Row 5 (was row 6): Hawthorn v Fremantle
$ws.Range("B5").Value = -7.75
$ws.Range("C5").Value = 7.75

# Row 6 (was row 7): North Melbourne v Gold Coast
$ws.Range("B6").Value = 24.5
$ws.Range("C6").Value = -24.5

# Row 7 (was row 8): Port Adelaide v Richmond
$ws.Range("B7").Value = -40
$ws.Range("C7").Value = 40

# Row 8 (was row 9): Sydney v Melbourne
$ws.Range("B8").Value = -2.75
$ws.Range("C8").Value = 2.75

# Row 9 (was row 10): West Coast v Adelaide
$ws.Range("B9").Value = 31.25
$ws.Range("C9").Value = -31.25
